# docs: corrected column name in group table in document. added insert
# statement for Group reference data table
#
# This script applies six localized corrections to the document:
#   1) Marks "favicon" with spell-check proof-error markers.
#   2) Marks every standalone "url" (in the "Each card should be a short
#      url..." paragraph) with spell-check proof-error markers.
#   3) Collapses the split "VARCHAR (" / "45" / ")" runs into one run.
#   4) Collapses the split "VARCHAR(" / "100" / ")" runs into one run.
#   5) Collapses the split "Modified " / "Operator" runs into one run.
#   6) Renames the GROUPTYPEVALUE column to GROUPVALUE in the Group table.
#
# Every edit is applied the same way: locate the target paragraph with
# Find, then replace that whole paragraph's Range via InsertXML with a
# freshly-built <w:p> (Flat-OPC wrapped) that has the exact desired run
# layout. Word always treats a "whole paragraph" InsertXML as a full
# paragraph content replacement, which lets us freely add <w:proofErr/>
# markers and merge/split runs with surgical precision.

$d = $word.ActiveDocument

$pkgOpen = "<pkg:package xmlns:pkg='http://schemas.microsoft.com/office/2006/xmlPackage'><pkg:part pkg:name='/word/document.xml' pkg:contentType='application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml'><pkg:xmlData><w:document xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'><w:body>"
$pkgClose = "</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>"

function Replace-Paragraph($para, [string]$innerXml) {
    $prng = $para.Range
    $prng.InsertXML($pkgOpen + $innerXml + $pkgClose)
}

# ---------------------------------------------------------------------
# 1) "... Default picture would be the favicon of the serving application."
# ---------------------------------------------------------------------
$rng = $d.Content
$rng.Find.Execute("Default picture would be the favicon", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$para = $rng.Paragraphs(1)

$p1 = "<w:p><w:pPr><w:spacing w:after=`"0`" w:line=`"240`" w:lineRule=`"auto`"/></w:pPr>" +
      "<w:r><w:t xml:space=`"preserve`">2. User should be able to create cards representing the </w:t></w:r>" +
      "<w:r><w:t>URL</w:t></w:r>" +
      "<w:r><w:t xml:space=`"preserve`"> where each card has a short title, brief description and a customizable picture. Default picture would be the </w:t></w:r>" +
      "<w:proofErr w:type=`"spellStart`"/>" +
      "<w:r><w:t>favicon</w:t></w:r>" +
      "<w:proofErr w:type=`"spellEnd`"/>" +
      "<w:r><w:t xml:space=`"preserve`"> of the serving application.</w:t></w:r>" +
      "<w:r><w:tab/></w:r>" +
      "</w:p>"
Replace-Paragraph $para $p1

# ---------------------------------------------------------------------
# 2) "Each card should be a short url ... contextual information too."
# ---------------------------------------------------------------------
$rng = $d.Content
$rng.Find.Execute("Each card should be a short url", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$para = $rng.Paragraphs(1)

$p2 = "<w:p><w:pPr><w:spacing w:after=`"0`" w:line=`"240`" w:lineRule=`"auto`"/></w:pPr>" +
      "<w:r><w:t xml:space=`"preserve`">4. </w:t></w:r>" +
      "<w:r><w:t xml:space=`"preserve`">Each card should be a short </w:t></w:r>" +
      "<w:proofErr w:type=`"spellStart`"/>" +
      "<w:r><w:t>url</w:t></w:r>" +
      "<w:proofErr w:type=`"spellEnd`"/>" +
      "<w:r><w:t xml:space=`"preserve`"> with the re-direction to the original </w:t></w:r>" +
      "<w:proofErr w:type=`"spellStart`"/>" +
      "<w:r><w:t>url</w:t></w:r>" +
      "<w:proofErr w:type=`"spellEnd`"/>" +
      "<w:r><w:t xml:space=`"preserve`">. </w:t></w:r>" +
      "<w:r><w:rPr><w:b/></w:rPr><w:t xml:space=`"preserve`">This short </w:t></w:r>" +
      "<w:proofErr w:type=`"spellStart`"/>" +
      "<w:r><w:rPr><w:b/></w:rPr><w:t>url</w:t></w:r>" +
      "<w:proofErr w:type=`"spellEnd`"/>" +
      "<w:r><w:rPr><w:b/></w:rPr><w:t xml:space=`"preserve`"> will have no expiration as it belongs to a group</w:t></w:r>" +
      "<w:r><w:t xml:space=`"preserve`">. The generation of this short </w:t></w:r>" +
      "<w:proofErr w:type=`"spellStart`"/>" +
      "<w:r><w:t>url</w:t></w:r>" +
      "<w:proofErr w:type=`"spellEnd`"/>" +
      "<w:r><w:t xml:space=`"preserve`"> should be dynamic and unique and could carry some contextual information too.</w:t></w:r>" +
      "</w:p>"
Replace-Paragraph $para $p2

# ---------------------------------------------------------------------
# Shared paragraph properties / run properties used by the DB-schema
# table cells (3-6).
# ---------------------------------------------------------------------
$tblRPr = "<w:rPr><w:rFonts w:asciiTheme=`"majorHAnsi`" w:hAnsiTheme=`"majorHAnsi`"/><w:sz w:val=`"20`"/></w:rPr>"

# ---------------------------------------------------------------------
# 3) VARCHAR (45)  -- merge "VARCHAR (" + "45" + ")" into one run
# ---------------------------------------------------------------------
$rng = $d.Content
$rng.Find.Execute("VARCHAR (45)", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$para = $rng.Paragraphs(1)

$p3 = "<w:p><w:pPr>$tblRPr</w:pPr><w:r>$tblRPr<w:t xml:space=`"preserve`">VARCHAR (45)</w:t></w:r></w:p>"
Replace-Paragraph $para $p3

# ---------------------------------------------------------------------
# 4) VARCHAR(100)  -- merge "VARCHAR(" + "100" + ")" into one run
# ---------------------------------------------------------------------
$rng = $d.Content
$rng.Find.Execute("VARCHAR(100)", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$para = $rng.Paragraphs(1)

$p4 = "<w:p><w:pPr>$tblRPr</w:pPr><w:r>$tblRPr<w:t xml:space=`"preserve`">VARCHAR(100)</w:t></w:r></w:p>"
Replace-Paragraph $para $p4

# ---------------------------------------------------------------------
# 5) "Modified Operator"  -- merge "Modified " + "Operator" into one run
# ---------------------------------------------------------------------
$rng = $d.Content
$rng.Find.Execute("Modified Operator", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$modifiedOperatorEnd = $rng.End
$para = $rng.Paragraphs(1)

$p5 = "<w:p><w:pPr>$tblRPr</w:pPr><w:r>$tblRPr<w:t xml:space=`"preserve`">Modified Operator</w:t></w:r></w:p>"
Replace-Paragraph $para $p5

# ---------------------------------------------------------------------
# 6) GROUPTYPEVALUE -> GROUPVALUE (only the 3rd occurrence, in the
#    Group reference table; search forward from the "Modified Operator"
#    row fixed above so the earlier two unrelated matches are skipped).
# ---------------------------------------------------------------------
$rng = $d.Content
$rng.Start = $modifiedOperatorEnd
$rng.Find.Execute("GROUPTYPEVALUE", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$para = $rng.Paragraphs(1)

$p6 = "<w:p><w:pPr>$tblRPr</w:pPr><w:r>$tblRPr<w:t xml:space=`"preserve`">GROUPVALUE</w:t></w:r></w:p>"
Replace-Paragraph $para $p6

Write-Output "edits applied"
